$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (avoid Excel auto-converting
# numeric-looking strings like "8.070" or "0.00000000357" into numbers,
# which would silently drop trailing zeros / switch to scientific notation).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.591.07"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.112.21"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.83%  "
$ws.Range("D5").Value = "348.64"
$ws.Range("E5").Value = "  +4.50%  "
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").Value = "0.5263"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "0.4519"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").Value = "53.69"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "0.09030"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "1.172"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "24.41"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "2.114.15"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "6.809"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "8.070"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "99.78"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("E17").Value = "  +4.38%  "
$ws.Range("D18").Value = "1.013"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "0.06702"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "19.34"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").Value = "6.322"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "30.639.21"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("E24").Value = "  +3.61%  "
$ws.Range("D25").Value = "2.390"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "2.354.90"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "22.29"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "165.26"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").Value = "2.525"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").Value = "135.19"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("D31").Value = "1.192"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "0.1072"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "1.641"
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").Value = "6.351"
$ws.Range("D35").Value = "4.003"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").Value = "5.916"
$ws.Range("E36").Value = "  +6.69%  "
$ws.Range("D37").Value = "10.19"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").Value = "0.02639"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("D39").Value = "0.06833"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "0.2319"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").Value = "12.61"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "0.6872"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "1.269"
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").Value = "14.85"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("D45").Value = "2.320"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "0.6416"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "3.759"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").Value = "0.00000000357"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "1.255"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "82.63"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "0.07288"
$ws.Range("E51").Value = "  +2.69%  "
